# Checklist workbook update:
#  - add a new checklist item about GET/POST + mysqldump/backup related note
#  - highlight row "localtime works correctly?" (A4) in green, like the other
#    already-finished/important items (A1, A3, A5, ...)
#  - turn wrap-text on for the whole checklist column so long notes are
#    readable instead of being clipped
#  - move the selection back up to A4 (top of sheet) instead of the bottom

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Append the new checklist entry in row 42
# ---------------------------------------------------------------------
$newRow = 42
$ws.Cells.Item($newRow, 1).Value = "when to use GET and POST in form when you send no data but this request updates the state of application (actually, not of application but of server file system. it creates some files, not related to db directly)"

# ---------------------------------------------------------------------
# 2) Highlight A4 ("localtime works correctly?") with the same green fill
#    used for the other highlighted rows (A1, A3, A5, A11, A15, A20, A24)
# ---------------------------------------------------------------------
$ws.Range("A4").Interior.Color = $ws.Range("A1").Interior.Color

# ---------------------------------------------------------------------
# 3) Turn on word-wrap for the whole used column A (rows 1-42), keeping
#    each cell's existing fill color intact
# ---------------------------------------------------------------------
for ($r = 1; $r -le $newRow; $r++) {
    $ws.Cells.Item($r, 1).WrapText = $true
}

# Make the new row tall enough to show the wrapped text
$ws.Rows.Item($newRow).RowHeight = 30

# ---------------------------------------------------------------------
# 4) Reset the view: no more frozen/top-left scroll to the bottom of the
#    sheet, select A4 near the top instead of A40
# ---------------------------------------------------------------------
$ws.Range("A4").Select() | Out-Null

Write-Output "checklist updated"
